$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 2270.9092
$ws.Range("I15").Value = 2270.9092
$ws.Range("K15").Value = 6812.7276
$ws.Range("M15").Value = -6643.7276

# Row 103: Let Loose the Juice
$ws.Range("H103").Value = 996.9
$ws.Range("I103").Value = 1210.5714
$ws.Range("J103").Value = 498.33334
$ws.Range("K103").Value = 3631.7142
$ws.Range("L103").Value = 1495.00002
$ws.Range("M103").Value = -3045.7142
$ws.Range("N103").Value = -2667.00002

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 7335.244
$ws.Range("I132").Value = 6380.1787
$ws.Range("J132").Value = 9392.308000000001
$ws.Range("K132").Value = 19140.5361
$ws.Range("L132").Value = 28176.924
$ws.Range("M132").Value = -16610.5361
$ws.Range("N132").Value = -33236.924

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 8773648
$ws.Range("I137").Value = 11112732
$ws.Range("J137").Value = 2083.3333
$ws.Range("K137").Value = 33338196
$ws.Range("L137").Value = 6249.999899999999
$ws.Range("M137").Value = -33335646
$ws.Range("N137").Value = -11349.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3480824.2
$ws.Range("I32").Value = 4790.1323
$ws.Range("J32").Value = 25673966
$ws.Range("K32").Value = 4790.1323
$ws.Range("L32").Value = 25673966
$ws.Range("M32").Value = -4503.1323
$ws.Range("N32").Value = -25674540

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2458.543
$ws.Range("I45").Value = 1404.3158
$ws.Range("J45").Value = 3710.4375
$ws.Range("K45").Value = 1404.3158
$ws.Range("L45").Value = 3710.4375
$ws.Range("M45").Value = -1027.3158
$ws.Range("N45").Value = -4464.4375

# Row 97: Ore for Me
$ws.Range("H97").Value = 248.42857
$ws.Range("I97").Value = 248.42857
$ws.Range("K97").Value = 248.42857
$ws.Range("M97").Value = 247.57143

# Row 107: Shielding the Realm
$ws.Range("H107").Value = 29980
$ws.Range("J107").Value = 29980
$ws.Range("L107").Value = 29980
$ws.Range("N107").Value = -37660

# Row 109: A Head of Demand
$ws.Range("H109").Value = 26667.334
$ws.Range("J109").Value = 26667.334
$ws.Range("L109").Value = 26667.334
$ws.Range("N109").Value = -29441.334

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 2143.8708
$ws.Range("I110").Value = 1353.826
$ws.Range("J110").Value = 4415.25
$ws.Range("K110").Value = 1353.826
$ws.Range("L110").Value = 4415.25
$ws.Range("M110").Value = 691.174
$ws.Range("N110").Value = -8505.25

# Row 112: Wrapped Knuckles
$ws.Range("H112").Value = 16296.75
$ws.Range("J112").Value = 16296.75
$ws.Range("L112").Value = 16296.75
$ws.Range("N112").Value = -19250.75

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2466.4
$ws.Range("I122").Value = 1211.5
$ws.Range("J122").Value = 3721.3
$ws.Range("K122").Value = 3634.5
$ws.Range("L122").Value = 11163.9
$ws.Range("M122").Value = -1184.5
$ws.Range("N122").Value = -16063.9

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1956446.5
$ws.Range("I132").Value = 1357.1428
$ws.Range("J132").Value = 5688890
$ws.Range("K132").Value = 4071.4284
$ws.Range("L132").Value = 17066670
$ws.Range("M132").Value = -1541.4284
$ws.Range("N132").Value = -17071730

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal
$ws.Range("H94").Value = 704.2143
$ws.Range("I94").Value = 714.4545000000001
$ws.Range("J94").Value = 666.6667
$ws.Range("K94").Value = 714.4545000000001
$ws.Range("L94").Value = 666.6667
$ws.Range("M94").Value = -263.4545000000001
$ws.Range("N94").Value = -1568.6667

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 2045
$ws.Range("I99").Value = 1504.6471
$ws.Range("K99").Value = 1504.6471
$ws.Range("M99").Value = -6.647099999999909

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2741.7827
$ws.Range("I105").Value = 1658
$ws.Range("J105").Value = 4773.875
$ws.Range("K105").Value = 1658
$ws.Range("L105").Value = 4773.875
$ws.Range("M105").Value = 89
$ws.Range("N105").Value = -8267.875

# Row 112: Enlistment Highs
$ws.Range("H112").Value = 35732.25
$ws.Range("J112").Value = 35732.25
$ws.Range("L112").Value = 35732.25
$ws.Range("N112").Value = -38686.25

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 3090.8667
$ws.Range("I16").Value = 3684.75
$ws.Range("J16").Value = 2412.1428
$ws.Range("K16").Value = 3684.75
$ws.Range("L16").Value = 2412.1428
$ws.Range("M16").Value = -3397.75
$ws.Range("N16").Value = -2986.1428

# Row 31: Wall Not Found
$ws.Range("H31").Value = 7814700.5
$ws.Range("I31").Value = 1649.7742
$ws.Range("J31").Value = 15154233
$ws.Range("K31").Value = 1649.7742
$ws.Range("L31").Value = 15154233
$ws.Range("M31").Value = -1354.7742
$ws.Range("N31").Value = -15154823

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 7814700.5
$ws.Range("I34").Value = 1649.7742
$ws.Range("J34").Value = 15154233
$ws.Range("K34").Value = 1649.7742
$ws.Range("L34").Value = 15154233
$ws.Range("M34").Value = -1447.7742
$ws.Range("N34").Value = -15154637

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 435857.7
$ws.Range("I105").Value = 626071.25
$ws.Range("J105").Value = 1083.8572
$ws.Range("K105").Value = 626071.25
$ws.Range("L105").Value = 1083.8572
$ws.Range("M105").Value = -624324.25
$ws.Range("N105").Value = -4577.8572

# Row 107: Built to Last
$ws.Range("H107").Value = 1180.9688
$ws.Range("I107").Value = 499.09525
$ws.Range("J107").Value = 2482.7273
$ws.Range("K107").Value = 499.09525
$ws.Range("L107").Value = 2482.7273
$ws.Range("M107").Value = 1420.90475
$ws.Range("N107").Value = -6322.7273

# Row 113: Patient Patients
$ws.Range("H113").Value = 3090.8667
$ws.Range("I113").Value = 3684.75
$ws.Range("J113").Value = 2412.1428
$ws.Range("K113").Value = 3684.75
$ws.Range("L113").Value = 2412.1428
$ws.Range("M113").Value = -1514.75
$ws.Range("N113").Value = -6752.1428

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 2104.628
$ws.Range("I134").Value = 1090.9231
$ws.Range("J134").Value = 3655
$ws.Range("K134").Value = 3272.7693
$ws.Range("L134").Value = 10965
$ws.Range("M134").Value = -737.7692999999999
$ws.Range("N134").Value = -16035

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Slippery Service
$ws.Range("H107").Value = 8803115
$ws.Range("I107").Value = 20834534
$ws.Range("J107").Value = 310348.75
$ws.Range("K107").Value = 62503602
$ws.Range("L107").Value = 931046.25
$ws.Range("M107").Value = -62501682
$ws.Range("N107").Value = -934886.25

# Row 132: More Mezcal
$ws.Range("H132").Value = 928.1539
$ws.Range("I132").Value = 767
$ws.Range("J132").Value = 976.5
$ws.Range("K132").Value = 6903
$ws.Range("L132").Value = 8788.5
$ws.Range("M132").Value = -4373
$ws.Range("N132").Value = -13848.5

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 1459670.9
$ws.Range("I102").Value = 2042282.6
$ws.Range("J102").Value = 3141.5
$ws.Range("K102").Value = 2042282.6
$ws.Range("L102").Value = 3141.5
$ws.Range("M102").Value = -2040660.6
$ws.Range("N102").Value = -6385.5

# Row 132: On Board for Lar
$ws.Range("H132").Value = 1159474.2
$ws.Range("I132").Value = 1985958.2
$ws.Range("J132").Value = 2396.5334
$ws.Range("K132").Value = 5957874.6
$ws.Range("L132").Value = 7189.600199999999
$ws.Range("M132").Value = -5955344.6
$ws.Range("N132").Value = -12249.6002

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 3156.4119
$ws.Range("I61").Value = 1241.7273
$ws.Range("J61").Value = 6666.6665
$ws.Range("K61").Value = 1241.7273
$ws.Range("L61").Value = 6666.6665
$ws.Range("M61").Value = -1039.7273
$ws.Range("N61").Value = -7070.6665

# Row 110: Breeches of Trust
$ws.Range("H110").Value = 22548
$ws.Range("J110").Value = 22548
$ws.Range("L110").Value = 22548
$ws.Range("N110").Value = -30728

# Row 113: Peace in Rest
$ws.Range("H113").Value = 3156.4119
$ws.Range("I113").Value = 1241.7273
$ws.Range("J113").Value = 6666.6665
$ws.Range("K113").Value = 1241.7273
$ws.Range("L113").Value = 6666.6665
$ws.Range("M113").Value = 928.2727
$ws.Range("N113").Value = -11006.6665

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 29445728
$ws.Range("I132").Value = 41712540
$ws.Range("J132").Value = 5378.8
$ws.Range("K132").Value = 125137620
$ws.Range("L132").Value = 16136.4
$ws.Range("M132").Value = -125135090
$ws.Range("N132").Value = -21196.4

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2078.8936
$ws.Range("J132").Value = 3659.8
$ws.Range("L132").Value = 10979.4
$ws.Range("N132").Value = -16039.4
